$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @{
    2 = @(44714, 80,  9000,  10000, 9500,  528)
    3 = @(45205, 200, 11000, 12000, 11500, 639)
    4 = @(44804, 50,  9500,  10000, 9750,  542)
    5 = @(45175, 250, 11000, 12000, 11500, 639)
    6 = @(45245, 100, 9000,  10000, 9500,  528)
    7 = @(45092, 210, 10000, 11000, 10714, 595)
    8 = @(45215, 200, 11000, 12000, 11500, 639)
    9 = @(44792, 160, 9000,  10000, 9500,  528)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value2  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value2 = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value2 = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value2 = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value2 = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value2 = $vals[5]   # P - Precio $/Kg
}
